$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$newTimestamps = @("2021-10-05 14:21:20.550749","2021-10-05 14:21:20.550758","2021-10-05 14:21:20.550761","2021-10-05 14:21:20.550764","2021-10-05 14:21:20.550767","2021-10-05 14:21:20.550770","2021-10-05 14:21:20.550773","2021-10-05 14:21:20.550775","2021-10-05 14:21:20.550778","2021-10-05 14:21:20.550781","2021-10-05 14:21:20.550784","2021-10-05 14:21:20.550786","2021-10-05 14:21:20.550789","2021-10-05 14:21:20.550791","2021-10-05 14:21:20.550794","2021-10-05 14:21:20.550797","2021-10-05 14:21:20.550800","2021-10-05 14:21:20.550803","2021-10-05 14:21:20.550805","2021-10-05 14:21:20.550808","2021-10-05 14:21:20.550811","2021-10-05 14:21:20.550814","2021-10-05 14:21:20.550816","2021-10-05 14:21:20.550819","2021-10-05 14:21:20.550822","2021-10-05 14:21:20.550824","2021-10-05 14:21:20.550827","2021-10-05 14:21:20.550830","2021-10-05 14:21:20.550832","2021-10-05 14:21:20.550835","2021-10-05 14:21:20.550838","2021-10-05 14:21:20.550840","2021-10-05 14:21:20.550843","2021-10-05 14:21:20.550846","2021-10-05 14:21:20.550848","2021-10-05 14:21:20.550851","2021-10-05 14:21:20.550854","2021-10-05 14:21:20.550856","2021-10-05 14:21:20.550859","2021-10-05 14:21:20.550862","2021-10-05 14:21:20.550865","2021-10-05 14:21:20.550868","2021-10-05 14:21:20.550870","2021-10-05 14:21:20.550873","2021-10-05 14:21:20.550876")
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $dataSheet.Cells.Item($i + 2, 6).Value = $newTimestamps[$i]
}

# --- Add the new "metadata" sheet, positioned after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Clone the "data" sheet's bold/bordered/centered header style (used on B1:F1
# there) onto this sheet's header row and the A2 index cell, so the new sheet
# matches the look of the existing one instead of inventing new style entries.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (row 1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Ketotic hypoglycaemia"
$metaSheet.Range("C2").Value = 248
$metaSheet.Range("D2").Value = "'1.5"
$metaSheet.Range("D2").ClearFormats()
$metaSheet.Range("E2").Value = "2021-08-26T10:52:37.646762Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:20.547353"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/248/?format=json"

$dataSheet.Activate()
